# handle array object and refactor createFiles method
#
# This script updates the "folder1" sheet (the key/en/ja/it translation
# table). The "a.b.list[1]" key is split into two keys
# ("a.b.list[1].x" and "a.b.list[1].y") to describe an array of objects,
# a new "a.b.list[2]" key is added, and a new "z" array (z[0], z[1], z[2])
# is appended at the bottom of the table. The "folder2" (hidden, language
# pivot) sheet content is unaffected by this change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("folder1")

# Make room: push the existing rows 6-10 ("a.d", "x.1", "x.2", "x.5", "y")
# down by two rows so we can fit the new "a.b.list[1].y" and
# "a.b.list[2]" rows right after "a.b.list[1]".
$ws.Rows("6:7").Insert()

# Row 5 used to be the single "a.b.list[1]" entry; it now becomes the
# ".x" member of that (now-array) object.
$ws.Range("A5").Value = "a.b.list[1].x"
$ws.Range("B5").Value = "a b list 1 x"
$ws.Range("C5").Value = "a b リスト 1 x"
$ws.Range("D5").Value = "a b elenco 1 x"

# New ".y" member of a.b.list[1].
$ws.Range("A6").Value = "a.b.list[1].y"
$ws.Range("B6").Value = "a b list 1 y"
$ws.Range("C6").Value = "a b リスト 1 y"
$ws.Range("D6").Value = "a b elenco 1 y"

# New a.b.list[2] entry.
$ws.Range("A7").Value = "a.b.list[2]"
$ws.Range("B7").Value = "a b list 2"
$ws.Range("C7").Value = "a b リスト 2"
$ws.Range("D7").Value = "a b elenco 2"

# Rows 8-12 already contain the shifted-down former rows 6-10
# (a.d, x.1, x.2, x.5, y) thanks to the Insert() above, so nothing to do
# for them.

# Append the new "z" array (z[0], z[1], z[2]) as rows 13-15.
$ws.Range("A13").Value = "z[0]"
$ws.Range("B13").Value = "z0"
$ws.Range("C13").Value = "z0 ja"
$ws.Range("D13").Value = "z0 it"

$ws.Range("A14").Value = "z[1]"
$ws.Range("B14").Value = "z1"
$ws.Range("C14").Value = "z1 ja"
$ws.Range("D14").Value = "z1 it"

$ws.Range("A15").Value = "z[2]"
$ws.Range("B15").Value = "z2"
$ws.Range("C15").Value = "z2 ja"
$ws.Range("D15").Value = "z2 it"

# Reflect the new data extent and restore the cursor position that was
# left selected on the sheet.
$ws.Range("B7").Select()
